$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1706.4286
$ws.Range("I106").Value = 1633.1072
$ws.Range("J106").Value = 1999.7142
$ws.Range("K106").Value = 1633.1072
$ws.Range("L106").Value = 1999.7142
$ws.Range("M106").Value = -1002.1072
$ws.Range("N106").Value = -3261.7142

$ws.Range("H112").Value = 1350.5
$ws.Range("J112").Value = 1364.6154
$ws.Range("L112").Value = 4093.8462
$ws.Range("N112").Value = -6309.8462

$ws.Range("H113").Value = 6149.143
$ws.Range("I113").Value = 1498.3334
$ws.Range("K113").Value = 1498.3334
$ws.Range("M113").Value = 1755.6666

$ws.Range("H116").Value = 6873.478
$ws.Range("I116").Value = 2108.6667
$ws.Range("J116").Value = 12071.454
$ws.Range("K116").Value = 2108.6667
$ws.Range("L116").Value = 12071.454
$ws.Range("M116").Value = 1333.3333
$ws.Range("N116").Value = -18955.454

$ws.Range("H132").Value = 24977826
$ws.Range("I132").Value = 28168506
$ws.Range("J132").Value = 2004933.2
$ws.Range("K132").Value = 84505518
$ws.Range("L132").Value = 6014799.6
$ws.Range("M132").Value = -84502988
$ws.Range("N132").Value = -6019859.6

$ws.Range("H137").Value = 713773.4399999999
$ws.Range("I137").Value = 1703988.9
$ws.Range("J137").Value = 2849.487
$ws.Range("K137").Value = 5111966.699999999
$ws.Range("L137").Value = 8548.460999999999
$ws.Range("M137").Value = -5109416.699999999
$ws.Range("N137").Value = -13648.461

$ws.Range("H138").Value = 3035.2263
$ws.Range("I138").Value = 1608.7894
$ws.Range("J138").Value = 3832.353
$ws.Range("K138").Value = 4826.3682
$ws.Range("L138").Value = 11497.059
$ws.Range("M138").Value = 313.6318000000001
$ws.Range("N138").Value = -21777.059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3825.3784
$ws.Range("I32").Value = 3674.862
$ws.Range("K32").Value = 3674.862
$ws.Range("M32").Value = -3387.862

$ws.Range("H61").Value = 2379.1428
$ws.Range("I61").Value = 2439.077
$ws.Range("J61").Value = 1600
$ws.Range("K61").Value = 2439.077
$ws.Range("L61").Value = 1600
$ws.Range("M61").Value = -2227.077
$ws.Range("N61").Value = -2024

$ws.Range("H74").Value = 281790.78
$ws.Range("I74").Value = 485093.78
$ws.Range("J74").Value = 2249.1875
$ws.Range("K74").Value = 485093.78
$ws.Range("L74").Value = 2249.1875
$ws.Range("M74").Value = -484219.78
$ws.Range("N74").Value = -3997.1875

$ws.Range("H77").Value = 281790.78
$ws.Range("I77").Value = 485093.78
$ws.Range("J77").Value = 2249.1875
$ws.Range("K77").Value = 2425468.9
$ws.Range("L77").Value = 11245.9375
$ws.Range("M77").Value = -2421100.9
$ws.Range("N77").Value = -19981.9375

$ws.Range("H102").Value = 1500
$ws.Range("I102").Value = 1500
$ws.Range("K102").Value = 1500
$ws.Range("M102").Value = 122

$ws.Range("H103").Value = 35000
$ws.Range("J103").Value = 35000
$ws.Range("L103").Value = 35000
$ws.Range("N103").Value = -37344

$ws.Range("H122").Value = 3510.7778
$ws.Range("I122").Value = 3284.8076
$ws.Range("J122").Value = 4098.3
$ws.Range("K122").Value = 9854.4228
$ws.Range("L122").Value = 12294.9
$ws.Range("M122").Value = -7404.4228
$ws.Range("N122").Value = -17194.9

$ws.Range("H132").Value = 3349.4614
$ws.Range("I132").Value = 2759.611
$ws.Range("J132").Value = 4676.625
$ws.Range("K132").Value = 8278.832999999999
$ws.Range("L132").Value = 14029.875
$ws.Range("M132").Value = -5748.832999999999
$ws.Range("N132").Value = -19089.875

$ws.Range("H136").Value = 2379.1428
$ws.Range("I136").Value = 2439.077
$ws.Range("J136").Value = 1600
$ws.Range("K136").Value = 7317.231000000001
$ws.Range("L136").Value = 4800
$ws.Range("M136").Value = -4767.231000000001
$ws.Range("N136").Value = -9900

$ws.Range("H137").Value = 38936
$ws.Range("J137").Value = 38936
$ws.Range("L137").Value = 38936
$ws.Range("N137").Value = -49136

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3462.2222
$ws.Range("I99").Value = 1556
$ws.Range("J99").Value = 5845
$ws.Range("K99").Value = 1556
$ws.Range("L99").Value = 5845
$ws.Range("M99").Value = -58
$ws.Range("N99").Value = -8841

$ws.Range("H103").Value = 36500
$ws.Range("J103").Value = 36500
$ws.Range("L103").Value = 36500
$ws.Range("N103").Value = -38844

$ws.Range("H134").Value = 3328.442
$ws.Range("I134").Value = 1123.24
$ws.Range("J134").Value = 6391.222
$ws.Range("K134").Value = 3369.72
$ws.Range("L134").Value = 19173.666
$ws.Range("M134").Value = -834.7200000000003
$ws.Range("N134").Value = -24243.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 197916.25
$ws.Range("I31").Value = 520478.34
$ws.Range("J31").Value = 2878.6978
$ws.Range("K31").Value = 520478.34
$ws.Range("L31").Value = 2878.6978
$ws.Range("M31").Value = -520183.34
$ws.Range("N31").Value = -3468.6978

$ws.Range("H34").Value = 197916.25
$ws.Range("I34").Value = 520478.34
$ws.Range("J34").Value = 2878.6978
$ws.Range("K34").Value = 520478.34
$ws.Range("L34").Value = 2878.6978
$ws.Range("M34").Value = -520276.34
$ws.Range("N34").Value = -3282.6978

$ws.Range("H99").Value = 4390.5
$ws.Range("J99").Value = 5773.778
$ws.Range("L99").Value = 5773.778
$ws.Range("N99").Value = -8769.778

$ws.Range("H103").Value = 34581
$ws.Range("I103").Value = 19662
$ws.Range("J103").Value = 49500
$ws.Range("K103").Value = 19662
$ws.Range("L103").Value = 49500
$ws.Range("M103").Value = -18490
$ws.Range("N103").Value = -51844

$ws.Range("H126").Value = 4390.5
$ws.Range("J126").Value = 5773.778
$ws.Range("L126").Value = 17321.334
$ws.Range("N126").Value = -22261.334

$ws.Range("H137").Value = 44940
$ws.Range("J137").Value = 44940
$ws.Range("L137").Value = 44940
$ws.Range("N137").Value = -55140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 495284.8
$ws.Range("I5").Value = 481.33334
$ws.Range("J5").Value = 636657.25
$ws.Range("K5").Value = 1444.00002
$ws.Range("L5").Value = 1909971.75
$ws.Range("M5").Value = -1332.00002
$ws.Range("N5").Value = -1910195.75

$ws.Range("H17").Value = 2350.2
$ws.Range("I17").Value = 1500.5
$ws.Range("J17").Value = 2916.6667
$ws.Range("K17").Value = 4501.5
$ws.Range("L17").Value = 8750.000100000001
$ws.Range("M17").Value = -4332.5
$ws.Range("N17").Value = -9088.000100000001

$ws.Range("H34").Value = 12182.4
$ws.Range("I34").Value = 16938.572
$ws.Range("J34").Value = 9621.385
$ws.Range("K34").Value = 50815.716
$ws.Range("L34").Value = 28864.155
$ws.Range("M34").Value = -50731.716
$ws.Range("N34").Value = -29032.155

$ws.Range("H39").Value = 21887.445
$ws.Range("J39").Value = 21887.445
$ws.Range("L39").Value = 65662.33499999999
$ws.Range("N39").Value = -66250.33499999999

$ws.Range("H51").Value = 1140
$ws.Range("I51").Value = 1050
$ws.Range("J51").Value = 1500
$ws.Range("K51").Value = 3150
$ws.Range("L51").Value = 4500
$ws.Range("M51").Value = -2690
$ws.Range("N51").Value = -5420

$ws.Range("H55").Value = 9302.5
$ws.Range("J55").Value = 9302.5
$ws.Range("L55").Value = 27907.5
$ws.Range("N55").Value = -28261.5

$ws.Range("H68").Value = 3608.4783
$ws.Range("I68").Value = 1139.4286
$ws.Range("J68").Value = 4688.6875
$ws.Range("K68").Value = 3418.2858
$ws.Range("L68").Value = 14066.0625
$ws.Range("M68").Value = -2607.2858
$ws.Range("N68").Value = -15688.0625

$ws.Range("H71").Value = 3608.4783
$ws.Range("I71").Value = 1139.4286
$ws.Range("J71").Value = 4688.6875
$ws.Range("K71").Value = 10254.8574
$ws.Range("L71").Value = 42198.1875
$ws.Range("M71").Value = -6198.857399999999
$ws.Range("N71").Value = -50310.1875

$ws.Range("H86").Value = 911.24
$ws.Range("J86").Value = 1310.9166
$ws.Range("L86").Value = 3932.7498
$ws.Range("N86").Value = -6304.7498

$ws.Range("H89").Value = 911.24
$ws.Range("J89").Value = 1310.9166
$ws.Range("L89").Value = 11798.2494
$ws.Range("N89").Value = -23654.2494

$ws.Range("H113").Value = 1603246.5
$ws.Range("I113").Value = 629.6909000000001
$ws.Range("J113").Value = 5435591
$ws.Range("K113").Value = 1889.0727
$ws.Range("L113").Value = 16306773
$ws.Range("M113").Value = 280.9272999999998
$ws.Range("N113").Value = -16311113

$ws.Range("H131").Value = 863.55
$ws.Range("J131").Value = 880.65265
$ws.Range("L131").Value = 2641.95795
$ws.Range("N131").Value = -12721.95795

$ws.Range("H132").Value = 2490.7058
$ws.Range("I132").Value = 1122.5
$ws.Range("J132").Value = 2673.1333
$ws.Range("K132").Value = 10102.5
$ws.Range("L132").Value = 24058.1997
$ws.Range("M132").Value = -7572.5
$ws.Range("N132").Value = -29118.1997

$ws.Range("H135").Value = 495284.8
$ws.Range("I135").Value = 481.33334
$ws.Range("J135").Value = 636657.25
$ws.Range("K135").Value = 4332.00006
$ws.Range("L135").Value = 5729915.25
$ws.Range("M135").Value = -1797.00006
$ws.Range("N135").Value = -5734985.25

$ws.Range("H140").Value = 2685.7896
$ws.Range("I140").Value = 730
$ws.Range("J140").Value = 3588.4614
$ws.Range("K140").Value = 2190
$ws.Range("L140").Value = 10765.3842
$ws.Range("M140").Value = 2990
$ws.Range("N140").Value = -21125.3842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 10285
$ws.Range("J109").Value = 10285
$ws.Range("L109").Value = 10285
$ws.Range("N109").Value = -12365

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3955.853
$ws.Range("I132").Value = 2814.4546
$ws.Range("J132").Value = 6048.4165
$ws.Range("K132").Value = 8443.363799999999
$ws.Range("L132").Value = 18145.2495
$ws.Range("M132").Value = -5913.363799999999
$ws.Range("N132").Value = -23205.2495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3206.7715
$ws.Range("I122").Value = 1825.5294
$ws.Range("J122").Value = 4511.278
$ws.Range("K122").Value = 5476.5882
$ws.Range("L122").Value = 13533.834
$ws.Range("M122").Value = -3026.5882
$ws.Range("N122").Value = -18433.834

$ws.Range("H132").Value = 2816.5
$ws.Range("I132").Value = 1149.9166
$ws.Range("J132").Value = 6149.6665
$ws.Range("K132").Value = 3449.7498
$ws.Range("L132").Value = 18448.9995
$ws.Range("M132").Value = -919.7498000000001
$ws.Range("N132").Value = -23508.9995
